$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eeg")

# Update the "Where will the EEG test take place?" paragraph (row 30) so the
# building name becomes a link (English + Dutch columns), per the commit's
# "psytoolkit link" / text-update changes.
$ws.Range("B30").Value = "EEG tests will occur at the <a href=\contact>Pieter de la Court building of Leiden University (Wassenaarseweg 52, 2333 AK Leiden)</a>. "
$ws.Range("C30").Value = "De EEG metingen zullen plaatsvinden in het  <a href=/contact>Pieder de la Court gebouw van Universiteit Leiden plaatsvinden (Wassenaarseweg 52, 2333 AK Leiden)</a>. "

# Restore the view state (scroll position + active selection) recorded when
# the workbook was last saved.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C30").Select()
